$wb = $excel.ActiveWorkbook

# Helper: set a cell's value as plain TEXT even when the text looks like a
# number/date, without Excel auto-converting it to a numeric/date type and
# without introducing a new cell style (keeps the default style index 0).
function Set-TextValue($cell, $text) {
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

# --- "Simple Fields" and "Simple Fields - Formatted" sheets ---
# Insert a new column D (tax-id) before the old billing-addr column; this
# shifts billing-addr (old D) -> E, invoice-no (old E) -> F, date (old F)
# -> G, total (old G) -> H, net-amount (old H) -> I, preserving their
# original (already-correct) values/types untouched. Then append a new
# trailing column J (items).
foreach ($sheetName in @("Simple Fields", "Simple Fields - Formatted")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Columns.Item(4).Insert()

    # New header cells
    $ws.Cells.Item(1, 4).Value = "tax-id"
    $ws.Cells.Item(1, 10).Value = "items"

    # billing-name (C2) and billing-addr (now E2) got re-OCR'd with new text
    $ws.Cells.Item(2, 3).Value = "บริษัท เอสซีจี เคมิคอลส์ จำกัด (สำนักงานใหญ่)"
    Set-TextValue $ws.Cells.Item(2, 4) "0105538052728"
    $ws.Cells.Item(2, 5).Value = "เลขที่ 1 ถ.ปูนซิเมนต์ไทย บางซื่อ กรุงเทพฯ 10800"
    $ws.Cells.Item(2, 10).Value = "table"
}

# --- "items" and "items - Formatted" sheets ---
# Insert two new columns (B and C) for quantity / unit-price; this shifts
# the old line-amount column (B) to D, preserving its value/type untouched.
foreach ($sheetName in @("items", "items - Formatted")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Columns.Item(2).Insert()
    $ws.Columns.Item(2).Insert()

    # New header cells
    $ws.Cells.Item(1, 2).Value = "quantity"
    $ws.Cells.Item(1, 3).Value = "unit-price"

    # description text was re-OCR'd slightly, quantity is a brand new value;
    # unit-price (C2) is left blank (no cell), matching the source invoice.
    Set-TextValue $ws.Cells.Item(2, 1) "50% Deposit Fee Produciton Course of Digital Learning 5 steps of IBE Competency Project Managemnet Fee Included"
    Set-TextValue $ws.Cells.Item(2, 2) "1"
}
